# Natmi following Dr Hou advice
# Update recomputed NATMI ligand-receptor metrics for Fbn1-Itgav (rows 2-17)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"3"
$ws.Range("G2").Value = [double]"11.98327633333333"
$ws.Range("H2").Value = [double]"35.949829"
$ws.Range("I2").Value = [double]"0.03345300399843466"
$ws.Range("J2").Value = [double]"0.03345300399843466"
$ws.Range("K2").Value = [double]"3"
$ws.Range("M2").Value = [double]"13.441269"
$ws.Range("N2").Value = [double]"40.323807"
$ws.Range("O2").Value = [double]"0.08973082133481231"
$ws.Range("P2").Value = [double]"0.08973082133481232"
$ws.Range("Q2").Value = [double]"161.070440697667"
$ws.Range("R2").Value = [double]"1449.633966279003"
$ws.Range("S2").Value = [double]"0.003001765524896302"
$ws.Range("T2").Value = [double]"0.003001765524896302"
$ws.Range("E3").Value = [double]"3"
$ws.Range("G3").Value = [double]"11.98327633333333"
$ws.Range("H3").Value = [double]"35.949829"
$ws.Range("I3").Value = [double]"0.03345300399843466"
$ws.Range("J3").Value = [double]"0.03345300399843466"
$ws.Range("K3").Value = [double]"3"
$ws.Range("M3").Value = [double]"54.711535"
$ws.Range("N3").Value = [double]"164.134605"
$ws.Range("O3").Value = [double]"0.3652416280068742"
$ws.Range("P3").Value = [double]"0.3652416280068742"
$ws.Range("Q3").Value = [double]"655.6234425258384"
$ws.Range("R3").Value = [double]"5900.610982732545"
$ws.Range("S3").Value = [double]"0.01221842964210874"
$ws.Range("T3").Value = [double]"0.01221842964210875"
$ws.Range("E4").Value = [double]"3"
$ws.Range("G4").Value = [double]"11.98327633333333"
$ws.Range("H4").Value = [double]"35.949829"
$ws.Range("I4").Value = [double]"0.03345300399843466"
$ws.Range("J4").Value = [double]"0.03345300399843466"
$ws.Range("K4").Value = [double]"3"
$ws.Range("M4").Value = [double]"63.67711"
$ws.Range("N4").Value = [double]"191.03133"
$ws.Range("O4").Value = [double]"0.4250937452800914"
$ws.Range("P4").Value = [double]"0.4250937452800915"
$ws.Range("Q4").Value = [double]"763.0604052380634"
$ws.Range("R4").Value = [double]"6867.54364714257"
$ws.Range("S4").Value = [double]"0.01422066276056446"
$ws.Range("T4").Value = [double]"0.01422066276056446"
$ws.Range("E5").Value = [double]"3"
$ws.Range("G5").Value = [double]"11.98327633333333"
$ws.Range("H5").Value = [double]"35.949829"
$ws.Range("I5").Value = [double]"0.03345300399843466"
$ws.Range("J5").Value = [double]"0.03345300399843466"
$ws.Range("K5").Value = [double]"3"
$ws.Range("M5").Value = [double]"17.96553866666667"
$ws.Range("N5").Value = [double]"53.896616"
$ws.Range("O5").Value = [double]"0.119933805378222"
$ws.Range("P5").Value = [double]"0.119933805378222"
$ws.Range("Q5").Value = [double]"215.2860143198516"
$ws.Range("R5").Value = [double]"1937.574128878664"
$ws.Range("S5").Value = [double]"0.004012146070865145"
$ws.Range("T5").Value = [double]"0.004012146070865145"
$ws.Range("E6").Value = [double]"3"
$ws.Range("G6").Value = [double]"332.3726806666667"
$ws.Range("H6").Value = [double]"997.1180420000001"
$ws.Range("I6").Value = [double]"0.9278651602470024"
$ws.Range("J6").Value = [double]"0.9278651602470025"
$ws.Range("K6").Value = [double]"3"
$ws.Range("M6").Value = [double]"13.441269"
$ws.Range("N6").Value = [double]"40.323807"
$ws.Range("O6").Value = [double]"0.08973082133481231"
$ws.Range("P6").Value = [double]"0.08973082133481232"
$ws.Range("Q6").Value = [double]"4467.510609091766"
$ws.Range("R6").Value = [double]"40207.5954818259"
$ws.Range("S6").Value = [double]"0.08325810291692076"
$ws.Range("T6").Value = [double]"0.08325810291692078"
$ws.Range("E7").Value = [double]"3"
$ws.Range("G7").Value = [double]"332.3726806666667"
$ws.Range("H7").Value = [double]"997.1180420000001"
$ws.Range("I7").Value = [double]"0.9278651602470024"
$ws.Range("J7").Value = [double]"0.9278651602470025"
$ws.Range("K7").Value = [double]"3"
$ws.Range("M7").Value = [double]"54.711535"
$ws.Range("N7").Value = [double]"164.134605"
$ws.Range("O7").Value = [double]"0.3652416280068742"
$ws.Range("P7").Value = [double]"0.3652416280068742"
$ws.Range("Q7").Value = [double]"18184.61955133816"
$ws.Range("R7").Value = [double]"163661.5759620434"
$ws.Range("S7").Value = [double]"0.3388949816994743"
$ws.Range("T7").Value = [double]"0.3388949816994744"
$ws.Range("E8").Value = [double]"3"
$ws.Range("G8").Value = [double]"332.3726806666667"
$ws.Range("H8").Value = [double]"997.1180420000001"
$ws.Range("I8").Value = [double]"0.9278651602470024"
$ws.Range("J8").Value = [double]"0.9278651602470025"
$ws.Range("K8").Value = [double]"3"
$ws.Range("M8").Value = [double]"63.67711"
$ws.Range("N8").Value = [double]"191.03133"
$ws.Range("O8").Value = [double]"0.4250937452800914"
$ws.Range("P8").Value = [double]"0.4250937452800915"
$ws.Range("Q8").Value = [double]"21164.53174780621"
$ws.Range("R8").Value = [double]"190480.7857302559"
$ws.Range("S8").Value = [double]"0.3944296760843105"
$ws.Range("T8").Value = [double]"0.3944296760843105"
$ws.Range("E9").Value = [double]"3"
$ws.Range("G9").Value = [double]"332.3726806666667"
$ws.Range("H9").Value = [double]"997.1180420000001"
$ws.Range("I9").Value = [double]"0.9278651602470024"
$ws.Range("J9").Value = [double]"0.9278651602470025"
$ws.Range("K9").Value = [double]"3"
$ws.Range("M9").Value = [double]"17.96553866666667"
$ws.Range("N9").Value = [double]"53.896616"
$ws.Range("O9").Value = [double]"0.119933805378222"
$ws.Range("P9").Value = [double]"0.119933805378222"
$ws.Range("Q9").Value = [double]"5971.254246260653"
$ws.Range("R9").Value = [double]"53741.28821634588"
$ws.Range("S9").Value = [double]"0.1112823995462967"
$ws.Range("T9").Value = [double]"0.1112823995462968"
$ws.Range("E10").Value = [double]"3"
$ws.Range("G10").Value = [double]"0.08615933333333332"
$ws.Range("H10").Value = [double]"0.258478"
$ws.Range("I10").Value = [double]"0.0002405259164795302"
$ws.Range("J10").Value = [double]"0.0002405259164795302"
$ws.Range("K10").Value = [double]"3"
$ws.Range("M10").Value = [double]"13.441269"
$ws.Range("N10").Value = [double]"40.323807"
$ws.Range("O10").Value = [double]"0.08973082133481231"
$ws.Range("P10").Value = [double]"0.08973082133481232"
$ws.Range("Q10").Value = [double]"1.158090776194"
$ws.Range("R10").Value = [double]"10.422816985746"
$ws.Range("S10").Value = [double]"2.158258803801671E-05"
$ws.Range("T10").Value = [double]"2.158258803801671E-05"
$ws.Range("E11").Value = [double]"3"
$ws.Range("G11").Value = [double]"0.08615933333333332"
$ws.Range("H11").Value = [double]"0.258478"
$ws.Range("I11").Value = [double]"0.0002405259164795302"
$ws.Range("J11").Value = [double]"0.0002405259164795302"
$ws.Range("K11").Value = [double]"3"
$ws.Range("M11").Value = [double]"54.711535"
$ws.Range("N11").Value = [double]"164.134605"
$ws.Range("O11").Value = [double]"0.3652416280068742"
$ws.Range("P11").Value = [double]"0.3652416280068742"
$ws.Range("Q11").Value = [double]"4.713909381243333"
$ws.Range("R11").Value = [double]"42.42518443119"
$ws.Range("S11").Value = [double]"8.785007731282904E-05"
$ws.Range("T11").Value = [double]"8.785007731282905E-05"
$ws.Range("E12").Value = [double]"3"
$ws.Range("G12").Value = [double]"0.08615933333333332"
$ws.Range("H12").Value = [double]"0.258478"
$ws.Range("I12").Value = [double]"0.0002405259164795302"
$ws.Range("J12").Value = [double]"0.0002405259164795302"
$ws.Range("K12").Value = [double]"3"
$ws.Range("M12").Value = [double]"63.67711"
$ws.Range("N12").Value = [double]"191.03133"
$ws.Range("O12").Value = [double]"0.4250937452800914"
$ws.Range("P12").Value = [double]"0.4250937452800915"
$ws.Range("Q12").Value = [double]"5.486377346193333"
$ws.Range("R12").Value = [double]"49.37739611574"
$ws.Range("S12").Value = [double]"0.0001022460626732099"
$ws.Range("T12").Value = [double]"0.00010224606267321"
$ws.Range("E13").Value = [double]"3"
$ws.Range("G13").Value = [double]"0.08615933333333332"
$ws.Range("H13").Value = [double]"0.258478"
$ws.Range("I13").Value = [double]"0.0002405259164795302"
$ws.Range("J13").Value = [double]"0.0002405259164795302"
$ws.Range("K13").Value = [double]"3"
$ws.Range("M13").Value = [double]"17.96553866666667"
$ws.Range("N13").Value = [double]"53.896616"
$ws.Range("O13").Value = [double]"0.119933805378222"
$ws.Range("P13").Value = [double]"0.119933805378222"
$ws.Range("Q13").Value = [double]"1.547898834494222"
$ws.Range("R13").Value = [double]"13.931089510448"
$ws.Range("S13").Value = [double]"2.884718845547445E-05"
$ws.Range("T13").Value = [double]"2.884718845547446E-05"
$ws.Range("E14").Value = [double]"3"
$ws.Range("G14").Value = [double]"13.77014866666667"
$ws.Range("H14").Value = [double]"41.310446"
$ws.Range("I14").Value = [double]"0.03844130983808348"
$ws.Range("J14").Value = [double]"0.03844130983808348"
$ws.Range("K14").Value = [double]"3"
$ws.Range("M14").Value = [double]"13.441269"
$ws.Range("N14").Value = [double]"40.323807"
$ws.Range("O14").Value = [double]"0.08973082133481231"
$ws.Range("P14").Value = [double]"0.08973082133481232"
$ws.Range("Q14").Value = [double]"185.088272398658"
$ws.Range("R14").Value = [double]"1665.794451587922"
$ws.Range("S14").Value = [double]"0.003449370304957231"
$ws.Range("T14").Value = [double]"0.003449370304957232"
$ws.Range("E15").Value = [double]"3"
$ws.Range("G15").Value = [double]"13.77014866666667"
$ws.Range("H15").Value = [double]"41.310446"
$ws.Range("I15").Value = [double]"0.03844130983808348"
$ws.Range("J15").Value = [double]"0.03844130983808348"
$ws.Range("K15").Value = [double]"3"
$ws.Range("M15").Value = [double]"54.711535"
$ws.Range("N15").Value = [double]"164.134605"
$ws.Range("O15").Value = [double]"0.3652416280068742"
$ws.Range("P15").Value = [double]"0.3652416280068742"
$ws.Range("Q15").Value = [double]"753.3859707315365"
$ws.Range("R15").Value = [double]"6780.473736583829"
$ws.Range("S15").Value = [double]"0.01404036658797828"
$ws.Range("T15").Value = [double]"0.01404036658797828"
$ws.Range("E16").Value = [double]"3"
$ws.Range("G16").Value = [double]"13.77014866666667"
$ws.Range("H16").Value = [double]"41.310446"
$ws.Range("I16").Value = [double]"0.03844130983808348"
$ws.Range("J16").Value = [double]"0.03844130983808348"
$ws.Range("K16").Value = [double]"3"
$ws.Range("M16").Value = [double]"63.67711"
$ws.Range("N16").Value = [double]"191.03133"
$ws.Range("O16").Value = [double]"0.4250937452800914"
$ws.Range("P16").Value = [double]"0.4250937452800915"
$ws.Range("Q16").Value = [double]"876.8432713636865"
$ws.Range("R16").Value = [double]"7891.58944227318"
$ws.Range("S16").Value = [double]"0.01634116037254333"
$ws.Range("T16").Value = [double]"0.01634116037254333"
$ws.Range("E17").Value = [double]"3"
$ws.Range("G17").Value = [double]"13.77014866666667"
$ws.Range("H17").Value = [double]"41.310446"
$ws.Range("I17").Value = [double]"0.03844130983808348"
$ws.Range("J17").Value = [double]"0.03844130983808348"
$ws.Range("K17").Value = [double]"3"
$ws.Range("M17").Value = [double]"17.96553866666667"
$ws.Range("N17").Value = [double]"53.896616"
$ws.Range("O17").Value = [double]"0.119933805378222"
$ws.Range("P17").Value = [double]"0.119933805378222"
$ws.Range("Q17").Value = [double]"247.3881383167484"
$ws.Range("R17").Value = [double]"2226.493244850736"
$ws.Range("S17").Value = [double]"0.004610412572604634"
$ws.Range("T17").Value = [double]"0.004610412572604636"

Write-Output "Applied 224 cell updates"
